$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the rate text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 2.9 = 10682.21 pesos`n✅ 10682.21 pesos = 2.89 = 969.81 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Sheet "tasas": update N10/O10 and N12/O12 ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 344.498
$ws2.Range("O10").Value = 3680
$ws2.Range("N12").Value = 3689.98
$ws2.Range("O12").Value = 335.002
